# Insert two new rows at position 196 (shifting existing rows 196:255 down to 198:257)
# and populate the two new rows with fresh data, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("196:197").Insert()

# Row 196: new "1a amarillo" record for Región de O'Higgins
$ws.Range("A196").Value = 11
$ws.Range("B196").Value = "Vega Monumental Concepción"
$ws.Range("C196").Value = "Bíobío"
$ws.Range("D196").Value = 44463
$ws.Range("E196").Value = 8
$ws.Range("F196").Value = "Fruta"
$ws.Range("G196").Value = 100102
$ws.Range("H196").Value = "Cítricos"
$ws.Range("I196").Value = 100102003
$ws.Range("J196").Value = "Limón"
$ws.Range("K196").Value = "Sin especificar"
$ws.Range("L196").Value = "1a amarillo"
$ws.Range("M196").Value = 300
$ws.Range("N196").Value = 5500
$ws.Range("O196").Value = 5500
$ws.Range("P196").Value = 5500
$ws.Range("Q196").Value = "`$/malla 16 kilos"
$ws.Range("R196").Value = "Región de O'Higgins"
$ws.Range("S196").Value = 344
$ws.Range("T196").Value = 16

# Row 197: new "2a amarillo" record for Región de O'Higgins
$ws.Range("A197").Value = 11
$ws.Range("B197").Value = "Vega Monumental Concepción"
$ws.Range("C197").Value = "Bíobío"
$ws.Range("D197").Value = 44463
$ws.Range("E197").Value = 8
$ws.Range("F197").Value = "Fruta"
$ws.Range("G197").Value = 100102
$ws.Range("H197").Value = "Cítricos"
$ws.Range("I197").Value = 100102003
$ws.Range("J197").Value = "Limón"
$ws.Range("K197").Value = "Sin especificar"
$ws.Range("L197").Value = "2a amarillo"
$ws.Range("M197").Value = 600
$ws.Range("N197").Value = 4000
$ws.Range("O197").Value = 4500
$ws.Range("P197").Value = 4250
$ws.Range("Q197").Value = "`$/malla 16 kilos"
$ws.Range("R197").Value = "Región de O'Higgins"
$ws.Range("S197").Value = 266
$ws.Range("T197").Value = 16
